# Update Clan Games data - 2026-02-22
# Sets column J (the "22/02/2026" Clan Games event column) to 0
# for every member row (rows 2 through 50) on the "clan games" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("clan games")

for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 10).Value = 0
}
